$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "59.735.38"
Set-TextValue $ws.Range("E2") "  +0.11%  "
Set-TextValue $ws.Range("D3") "2.530.29"
Set-TextValue $ws.Range("E3") "  +1.52%  "
Set-TextValue $ws.Range("E4") "  -0.23%  "
Set-TextValue $ws.Range("D5") "543.92"
Set-TextValue $ws.Range("E5") "  +0.11%  "
Set-TextValue $ws.Range("D6") "146.40"
Set-TextValue $ws.Range("E6") "  -0.50%  "
Set-TextValue $ws.Range("E7") "  -0.30%  "
Set-TextValue $ws.Range("E8") "  -1.06%  "
Set-TextValue $ws.Range("D9") "2.555.63"
Set-TextValue $ws.Range("E9") "  +1.36%  "
Set-TextValue $ws.Range("E10") "  +0.34%  "
Set-TextValue $ws.Range("E11") "  +0.23%  "
Set-TextValue $ws.Range("D12") "5.59"
Set-TextValue $ws.Range("E12") "  +2.51%  "
Set-TextValue $ws.Range("D13") "0.362"
Set-TextValue $ws.Range("E13") "  +1.32%  "
Set-TextValue $ws.Range("D14") "2.977.10"
Set-TextValue $ws.Range("E14") "  +0.18%  "
Set-TextValue $ws.Range("D15") "23.70"
Set-TextValue $ws.Range("E15") "  -3.72%  "
Set-TextValue $ws.Range("D16") "59.665.35"
Set-TextValue $ws.Range("E16") "  -0.27%  "
Set-TextValue $ws.Range("E17") "  +2.04%  "
Set-TextValue $ws.Range("D18") "2.539.26"
Set-TextValue $ws.Range("E18") "  +0.99%  "
Set-TextValue $ws.Range("E19") "  -1.66%  "
Set-TextValue $ws.Range("E20") "  -1.32%  "
Set-TextValue $ws.Range("D21") "327.56"
Set-TextValue $ws.Range("E21") "  +0.01%  "
Set-TextValue $ws.Range("E23") "  +2.33%  "
Set-TextValue $ws.Range("D24") "62.39"
Set-TextValue $ws.Range("E24") "  +1.45%  "
Set-TextValue $ws.Range("D25") "0.440"
Set-TextValue $ws.Range("E25") "  -2.07%  "
Set-TextValue $ws.Range("E26") "  +1.99%  "
Set-TextValue $ws.Range("E27") "  -1.51%  "
Set-TextValue $ws.Range("E28") "  +2.51%  "
Set-TextValue $ws.Range("D29") "0.0₃0801"
Set-TextValue $ws.Range("E29") "  +0.77%  "
Set-TextValue $ws.Range("D30") "6.90"
Set-TextValue $ws.Range("E30") "  +0.35%  "
Set-TextValue $ws.Range("E31") "  -0.05%  "
Set-TextValue $ws.Range("D32") "1.22"
Set-TextValue $ws.Range("E32") "  -7.02%  "
Set-TextValue $ws.Range("D33") "1.49"
Set-TextValue $ws.Range("E33") "  +2.93%  "
Set-TextValue $ws.Range("D34") "160.90"
Set-TextValue $ws.Range("E34") "  +1.27%  "
Set-TextValue $ws.Range("E35") "  +0.07%  "
Set-TextValue $ws.Range("D36") "18.81"
Set-TextValue $ws.Range("E36") "  -0.64%  "
Set-TextValue $ws.Range("D37") "4.45"
Set-TextValue $ws.Range("E37") "  -1.79%  "
Set-TextValue $ws.Range("E38") "  -6.66%  "
Set-TextValue $ws.Range("D39") "5.70"
Set-TextValue $ws.Range("E39") "  -5.85%  "
Set-TextValue $ws.Range("D40") "37.15"
Set-TextValue $ws.Range("E40") "  +0.90%  "
Set-TextValue $ws.Range("E41") "  +1.61%  "
Set-TextValue $ws.Range("D42") "301.23"
Set-TextValue $ws.Range("E42") "  -4.36%  "
Set-TextValue $ws.Range("E43") "  -1.68%  "
Set-TextValue $ws.Range("D45") "0.993"
Set-TextValue $ws.Range("E45") "  -0.14%  "
Set-TextValue $ws.Range("D46") "10.80"
Set-TextValue $ws.Range("E46") "  +0.08%  "
Set-TextValue $ws.Range("D47") "19.03"
Set-TextValue $ws.Range("E47") "  +1.66%  "
Set-TextValue $ws.Range("E48") "  -0.55%  "
Set-TextValue $ws.Range("D49") "123.70"
Set-TextValue $ws.Range("E49") "  -2.86%  "
Set-TextValue $ws.Range("D50") "0.0519"
Set-TextValue $ws.Range("E50") "  -2.72%  "
Set-TextValue $ws.Range("E51") "  -1.32%  "
